# Applies the "squid squad suggestions" content edits to the report-card
# table template (docs/table_template.xlsx).
#
# Semantic changes on Sheet1 (row/column positions are unchanged; only the
# text of a few cells was revised):
#   A3 - indicator label for "Number of commercial vessels" clarified
#   C2 - implication text for "Commercial landings" rewritten
#   C3 - implication text for "Number of commercial vessels" extended
#   C6 - implication text for "Bottom temperature" reworded/reordered
# Finally, the active selection moves from D4 to C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Order matches the author's edit order (new shared-string table entries are
# appended in first-use order, so this keeps the saved string indices
# aligned with the canonical edit).
$ws.Range("C2").Value = "Environmental dynamics vary between locations/timing of the summer and winter squid fisheries. An increase in landings since 2020 but decrease in number of vessels could indicate targeted trips in specific times of year and fishers targeting other species when longfin are not available. "

$ws.Range("A3").Value = "Number of commercial vessels (#  of federally-permitted vessels landing over 1lb of longfin squid)"

$ws.Range("C3").Value = "Number of commercial vessels has been steadily decreasing since around 2000 consistent with decreasing fleet diversity and continued risk to fishery resilience (MAFMC FID). Permit requalification in 2019 and a decrease in the incidental limit for trimester 2 resulted in fishery closures in 2022 and 2023, which may contribute to decreased participation."

$ws.Range("C6").Value = "Inshore temperature thresholds (around 14°C) initiate migration of squid from offshore overwintering habitats. Longfin squid seasonal distribution and growth rates are likely temperature dependent, avoiding water <8°C. 2024 spring bottom temperatures are near the long term mean, however cold pool temperatures dipped below 8°C. "

# Move the active selection to C6, matching the saved sheet view state.
$ws.Range("C6").Select()
